$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.996515393257141
$ws.Range("B1").Value = 2.262944936752319
$ws.Range("C1").Value = 2.239951372146606
$ws.Range("D1").Value = 2.680375576019287
$ws.Range("E1").Value = 1.520605802536011
